# Update device subsytems, components, equipment, ports and vessel dataspace.
#
# The "vessel_sf" sheet listed two vessel parameters that are no longer part
# of the safety-factor dataset: "Turntable outer diameter [m]" (row id 7)
# and "Dredge depth [m]" (row id 8). Remove those two rows; Excel re-flows
# the remaining rows upward (their original id numbers, 9-12, are kept as
# the values already stored in column A).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("vessel_sf")

# Make this the active sheet/tab, matching the saved selection state.
$ws.Activate()

# Rows 9 and 10 hold "Turntable outer diameter [m]" and "Dredge depth [m]".
$ws.Rows("9:10").Delete()

# Leave the selection where the deleted rows used to be.
$ws.Range("A9:XFD10").Select()
